$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.791.05'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.572.61'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.13'
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.84'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").Value = '2.576.36'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("E12").Value = '  +8.24%  '
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '3.026.63'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").Value = '58.885.35'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.06'
$ws.Range("E16").Value = '  +5.77%  '
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").Value = '2.574.38'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.64'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.92'
$ws.Range("E24").Value = '  -4.03%  '
$ws.Range("E25").Value = '  +6.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  +1.76%  '
$ws.Range("E28").Value = '  +1.25%  '
$ws.Range("D29").Value = '0.0₃0773'
$ws.Range("E29").Value = '  +1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.03'
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.95'
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.95'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.13'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.868'
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.872'
$ws.Range("E38").Value = '  +5.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.74'
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '289.89'
$ws.Range("E41").Value = '  +2.59%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.62'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0967'
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.61'
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.03'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.39'
$ws.Range("E49").Value = '  +8.29%  '
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.45'
$ws.Range("E51").Value = '  +3.05%  '
